# All enemies die when shot, music added
#
# 1) Turn the bare "BlackMetal.jpg" reference URL into a real hyperlink.
# 2) Append seven new bulleted reference entries (asteroid models/textures,
#    grate, beam texture, and three BGM credits) after it, preserving the
#    ListParagraph/numPr formatting and the proofErr spell-check wrappers
#    that Word leaves around the unusual proper nouns.
# 3) Re-home the "_GoBack" bookmark inside the final BGM paragraph, split
#    between "Carpenter Br" and "ut", exactly as the tracked change shows.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: convert the plain-text wallpoper.com URL into a hyperlink run.
# ---------------------------------------------------------------------
$wallpoperUrl = "http://wallpoper.com/images/00/24/09/44/metal-textures_00240944.jpg"
$fr = $d.Content
$fr.Find.Execute($wallpoperUrl) | Out-Null
$d.Hyperlinks.Add($fr, $wallpoperUrl) | Out-Null

# ---------------------------------------------------------------------
# Step 2: append the new list paragraphs (plain runs + proofErr markers)
# right before the end of the last paragraph's mark, so the existing
# bookmark / paragraph formatting on "BlackMetal..." is left untouched.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.End - 1, $lastPara.Range.End - 1)

$newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Asteroid models &#8211; From </w:t></w:r>
            <w:r><w:t>https://www.turbosquid.com/FullPreview/Index.cfm/ID/1051001</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Asteroid textures &#8211; From </w:t></w:r>
            <w:r><w:t>http://www.hedfiles.net/Ast_Rock_01.png</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t>Grate.png</w:t></w:r>
            <w:r><w:t xml:space="preserve"> &#8211; From </w:t></w:r>
            <w:r><w:t>http://www.vibrantmetaltech.com/img/e%202.png</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Beam texture &#8211;Edit of </w:t></w:r>
            <w:r><w:t>http://vortex-x.deviantart.com/art/Metal-Beam-118995632</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Lvl1 BGM &#8211; </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Technoir</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve">, by </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Perturbator</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Lvl2 BGM &#8211; Humans are Such Easy Prey, by </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Pertubator</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t>Lvl3 BGM &#8211; Turbo killer, by Carpenter Brut</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newParasXml) | Out-Null

# ---------------------------------------------------------------------
# Step 3: turn the newly-inserted plain-text URLs into real hyperlinks.
# ---------------------------------------------------------------------
$newLinks = @(
    "https://www.turbosquid.com/FullPreview/Index.cfm/ID/1051001",
    "http://www.hedfiles.net/Ast_Rock_01.png",
    "http://www.vibrantmetaltech.com/img/e%202.png",
    "http://vortex-x.deviantart.com/art/Metal-Beam-118995632"
)
foreach ($url in $newLinks) {
    $fr2 = $d.Content
    $fr2.Find.Execute($url) | Out-Null
    $d.Hyperlinks.Add($fr2, $url) | Out-Null
}

# ---------------------------------------------------------------------
# Step 4: move the "_GoBack" bookmark so that it once again sits right
# at the end of the document, splitting "Carpenter Br" | "ut" the same
# way Word leaves it after the author's last edit landed there.
# ---------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()
$splitFind = $d.Content
$splitFind.Find.Execute("Carpenter Br") | Out-Null
$splitPoint = $d.Range($splitFind.End, $splitFind.End)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
